# Updates cryptos list figures (price + 1h volume change) to the latest
# scraped values, and swaps the Bittensor/Stacks rows (rank order change).
#
# All cells in columns B-E are stored as plain text in this sheet, but a
# naive `.Value = "..."` assignment lets Excel's type-inference turn
# number-looking strings (e.g. "1.00", "600.82") into real numbers. Force
# text storage by flipping the cell to a text NumberFormat before the
# write, then snap the style back to "Normal" afterwards so no stray
# cell-style survives the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# r: row, d: new Price (column D), e: new Volume(1h) (column E)
$updates = @(
    @{ r = 2;  d = "68.998.62";  e = "  +0.42%  " },
    @{ r = 3;  d = "3.747.38";   e = "  +0.28%  " },
    @{ r = 4;  d = "1.00";       e = "  +0.01%  " },
    @{ r = 5;  d = "600.82";     e = "  -0.09%  " },
    @{ r = 6;  d = "166.10";     e = "  -1.69%  " },
    @{ r = 7;  d = "3.743.69";   e = "  +0.16%  " },
    @{ r = 9;  e = "  +0.86%  " },
    @{ r = 10; d = "0.173";      e = "  +5.59%  " },
    @{ r = 11; e = "  +0.99%  " },
    @{ r = 12; d = "0.460";      e = "  -0.37%  " },
    @{ r = 13; d = "37.78";      e = "  -1.03%  " },
    @{ r = 14; e = "  +1.33%  " },
    @{ r = 15; d = "4.373.85";   e = "  +0.22%  " },
    @{ r = 16; d = "3.752.60";   e = "  +0.37%  " },
    @{ r = 17; d = "69.117.72";  e = "  +0.55%  " },
    @{ r = 18; d = "7.44" },
    @{ r = 19; d = "17.66";      e = "  +2.99%  " },
    @{ r = 20; e = "  -0.90%  " },
    @{ r = 21; d = "11.28";      e = "  +5.37%  " },
    @{ r = 22; d = "490.76";     e = "  -0.65%  " },
    @{ r = 23; e = "  -0.40%  " },
    @{ r = 24; e = "  +3.54%  " },
    @{ r = 25; d = "84.58";      e = "  -0.85%  " },
    @{ r = 26; e = "  -1.81%  " },
    @{ r = 27; e = "  -0.72%  " },
    @{ r = 28; e = "  -0.75%  " },
    @{ r = 29; e = "  -0.05%  " },
    @{ r = 30; d = "2.97";       e = "  -0.26%  " },
    @{ r = 31; d = "8.17";       e = "  +3.35%  " },
    @{ r = 32; e = "  -4.43%  " },
    @{ r = 33; d = "31.71";      e = "  +0.10%  " },
    @{ r = 34; d = "3.888.98";   e = "  +0.11%  " },
    @{ r = 35; d = "3.682.54";   e = "  +0.26%  " },
    @{ r = 36; e = "  -0.19%  " },
    @{ r = 37; e = "  +2.05%  " },
    @{ r = 38; e = "  -0.01%  " },
    @{ r = 39; e = "  +3.88%  " },
    @{ r = 40; e = "  -0.11%  " },
    @{ r = 41; e = "  +9.46%  " },
    @{ r = 42; e = "  -0.02%  " },
    @{ r = 44; d = "48.65";      e = "  -0.57%  " },
    @{ r = 46; d = "8.45";       e = "  -0.68%  " },
    @{ r = 48; d = "40.08";      e = "  -1.33%  " },
    @{ r = 49; e = "  +10.39%  " },
    @{ r = 50; d = "141.34";     e = "  +0.12%  " },
    @{ r = 51; d = "2.793.00";   e = "  -0.10%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("d")) {
        Set-TextValue $ws.Cells.Item($u.r, 4) $u.d
    }
    if ($u.ContainsKey("e")) {
        Set-TextValue $ws.Cells.Item($u.r, 5) $u.e
    }
}

# Row 43 and 45 swap identity (Bittensor <-> Stacks) with refreshed values.
Set-TextValue $ws.Cells.Item(43, 2) "Stacks"
Set-TextValue $ws.Cells.Item(43, 3) "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Cells.Item(43, 4) "2.00"
Set-TextValue $ws.Cells.Item(43, 5) "  +0.69%  "

Set-TextValue $ws.Cells.Item(45, 2) "Bittensor"
Set-TextValue $ws.Cells.Item(45, 3) "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Cells.Item(45, 4) "427.21"
Set-TextValue $ws.Cells.Item(45, 5) "  -2.20%  "
